$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 110
$ws.Range("I2").Value = 273
$ws.Range("J2").Value = 1140
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 337
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = 206
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 12
$ws.Range("S2").Value = 126
$ws.Range("U2").Value = 12
$ws.Range("V2").Value = 1783
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1793
$ws.Range("Z2").Value = 33
